$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert 5 new rows before row 23 (pushes "Componentes a 1V2" section, which was
# at row 26, down to row 31) - formatting is inherited from the surrounding rows.
$ws.Rows("23:27").Insert()

# Fill in the newly inserted rows with the new "Componentes a 3V3" entries.
$ws.Range("A23").Value = "FPGA ICE"
$ws.Range("A24").Value = "Memoria NORFlash 32Mb"
$ws.Range("A25").Value = "Regulador lineal ajustable"

$ws.Range("D23").Value = "Board 1"
$ws.Range("D24").Value = "Board 1"
$ws.Range("D25").Value = "Board 1"

# Leave rows 26-30 blank, matching the original pattern.

# Update the active selection to reflect where the user last worked.
[void]$ws.Range("D28").Select()
